$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update I16: days-done count for week row 16 goes from 4 to 5.
$ws.Range("I16").Value = 5

# H16 picks up the same fill/border formatting as the rest of the row (style 17)
# by copying the format from a neighboring cell that already carries it.
$ws.Range("G16").Copy()
$ws.Range("H16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the final selection left by the editor.
$ws.Range("L16").Select()
